$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17.99061042797839
$ws.Range("C2").Value = 9.736186189957387
$ws.Range("D2").Value = 7.405907095948355
$ws.Range("F2").Value = 39.5819224932138
$ws.Range("G2").Value = 3.694501517802035
$ws.Range("I2").Value = 31.5548842137092
$ws.Range("L2").Value = 10.67039657122419
$ws.Range("M2").Value = 16.58805167610401
$ws.Range("B3").Value = 17.63132628203301
$ws.Range("C3").Value = 9.110936423818403
$ws.Range("D3").Value = 7.418924918032643
$ws.Range("F3").Value = 39.19784056479596
$ws.Range("G3").Value = 3.69860441482674
$ws.Range("I3").Value = 31.44829799129273
$ws.Range("L3").Value = 10.68345813809945
$ws.Range("M3").Value = 16.54133220765944
$ws.Range("B4").Value = 17.41311715636152
$ws.Range("C4").Value = 8.703597655823339
$ws.Range("D4").Value = 7.427654552514705
$ws.Range("F4").Value = 38.97089404488143
$ws.Range("G4").Value = 3.701252778963408
$ws.Range("I4").Value = 31.3893207206359
$ws.Range("L4").Value = 10.69293559898511
$ws.Range("M4").Value = 16.51649491236203
$ws.Range("B5").Value = 17.3249380912756
$ws.Range("C5").Value = 8.531703609648279
$ws.Range("D5").Value = 7.431396400340343
$ws.Range("F5").Value = 38.88072487710117
$ws.Range("G5").Value = 3.70236461862738
$ws.Range("I5").Value = 31.36691948606756
$ws.Range("L5").Value = 10.69716417969981
$ws.Range("M5").Value = 16.50734840591493
$ws.Range("B6").Value = 17.31034503825207
$ws.Range("C6").Value = 8.502804202461817
$ws.Range("D6").Value = 7.432028848815533
$ws.Range("F6").Value = 38.86589431271032
$ws.Range("G6").Value = 3.702551212099488
$ws.Range("I6").Value = 31.36329852138363
$ws.Range("L6").Value = 10.69788846082021
$ws.Range("M6").Value = 16.50588870487701
$ws.Range("B7").Value = 17.41192474743401
$ws.Range("C7").Value = 8.701303330594126
$ws.Range("D7").Value = 7.427704270605743
$ws.Range("F7").Value = 38.96966852357533
$ws.Range("G7").Value = 3.701267641416941
$ws.Range("I7").Value = 31.38901199323647
$ws.Range("L7").Value = 10.69299114355957
$ws.Range("M7").Value = 16.51636760347276
$ws.Range("B8").Value = 17.866317004443
$ws.Range("C8").Value = 9.525458033553328
$ws.Range("D8").Value = 7.410242280678305
$ws.Range("F8").Value = 39.44769737074365
$ws.Range("G8").Value = 3.695889465695571
$ws.Range("I8").Value = 31.51679536652118
$ws.Range("L8").Value = 10.67459763161124
$ws.Range("M8").Value = 16.57114858571871
$ws.Range("B9").Value = 18.77024427204812
$ws.Range("C9").Value = 10.95617882675835
$ws.Range("D9").Value = 7.381879365815343
$ws.Range("F9").Value = 40.451820455884
$ws.Range("G9").Value = 3.686361877666633
$ws.Range("I9").Value = 31.81833505333871
$ws.Range("L9").Value = 10.65009671745467
$ws.Range("M9").Value = 16.70875455574787
$ws.Range("B10").Value = 19.43423138102909
$ws.Range("C10").Value = 11.89509230529709
$ws.Range("D10").Value = 7.364673973617681
$ws.Range("F10").Value = 41.22489333074297
$ws.Range("G10").Value = 3.679974908001478
$ws.Range("I10").Value = 32.07034872067238
$ws.Range("L10").Value = 10.6391518825724
$ws.Range("M10").Value = 16.82771283879358
$ws.Range("B11").Value = 19.73467416903918
$ws.Range("C11").Value = 12.29805868762137
$ws.Range("D11").Value = 7.357646063046038
$ws.Range("F11").Value = 41.58305948832939
$ws.Range("G11").Value = 3.677200638775674
$ws.Range("I11").Value = 32.19144875328072
$ws.Range("L11").Value = 10.63570474913792
$ws.Range("M11").Value = 16.88557351621091
$ws.Range("B12").Value = 19.84808819005716
$ws.Range("C12").Value = 12.44719552162466
$ws.Range("D12").Value = 7.355100530084173
$ws.Range("F12").Value = 41.71951736499603
$ws.Range("G12").Value = 3.676168824652074
$ws.Range("I12").Value = 32.23821791050844
$ws.Range("L12").Value = 10.63461953698478
$ws.Range("M12").Value = 16.90800974491433
$ws.Range("B13").Value = 19.8236802651133
$ws.Range("C13").Value = 12.41522977695222
$ws.Range("D13").Value = 7.355643592205648
$ws.Range("F13").Value = 41.69009357608308
$ws.Range("G13").Value = 3.676390212730437
$ws.Range("I13").Value = 32.22810510466061
$ws.Range("L13").Value = 10.63484346842182
$ws.Range("M13").Value = 16.90315453933651
$ws.Range("B14").Value = 19.74401266261793
$ws.Range("C14").Value = 12.31039747592347
$ws.Range("D14").Value = 7.357434314848328
$ws.Range("F14").Value = 41.59426996999751
$ws.Range("G14").Value = 3.677115375911724
$ws.Range("I14").Value = 32.19527832382004
$ws.Range("L14").Value = 10.6356110569735
$ws.Range("M14").Value = 16.8874089009988
$ws.Range("B15").Value = 19.69516382845474
$ws.Range("C15").Value = 12.24573481369464
$ws.Range("D15").Value = 7.358546289536744
$ws.Range("F15").Value = 41.53567994703744
$ws.Range("G15").Value = 3.677561996301446
$ws.Range("I15").Value = 32.17528909279972
$ws.Range("L15").Value = 10.63610989171203
$ws.Range("M15").Value = 16.87783229710905
$ws.Range("B16").Value = 19.41455345794681
$ws.Range("C16").Value = 11.86827280666061
$ws.Range("D16").Value = 7.365149417216246
$ws.Range("F16").Value = 41.20160797599189
$ws.Range("G16").Value = 3.680158842370889
$ws.Range("I16").Value = 32.06256304521436
$ws.Range("L16").Value = 10.63940797430399
$ws.Range("M16").Value = 16.82400579179253
$ws.Range("B17").Value = 19.2419044343693
$ws.Range("C17").Value = 11.63053710155244
$ws.Range("D17").Value = 7.369405482496541
$ws.Range("F17").Value = 40.9982536982896
$ws.Range("G17").Value = 3.681785437031948
$ws.Range("I17").Value = 31.99505140129386
$ws.Range("L17").Value = 10.64182348694252
$ws.Range("M17").Value = 16.79193540679019
$ws.Range("B18").Value = 19.14245701321621
$ws.Range("C18").Value = 11.49152374791765
$ws.Range("D18").Value = 7.371928601589181
$ws.Range("F18").Value = 40.88190710299179
$ws.Range("G18").Value = 3.682733367395256
$ws.Range("I18").Value = 31.9568294763464
$ws.Range("L18").Value = 10.64335700993533
$ws.Range("M18").Value = 16.77384279059712
$ws.Range("B19").Value = 19.1087649947202
$ws.Range("C19").Value = 11.4440651980084
$ws.Range("D19").Value = 7.372795764787194
$ws.Range("F19").Value = 40.84262337502601
$ws.Range("G19").Value = 3.683056446297649
$ws.Range("I19").Value = 31.94399328610236
$ws.Range("L19").Value = 10.64390100087938
$ws.Range("M19").Value = 16.76777802499736
$ws.Range("B20").Value = 19.2602990828007
$ws.Range("C20").Value = 11.65607972389272
$ws.Range("D20").Value = 7.368944633556533
$ws.Range("F20").Value = 41.0198379685718
$ws.Range("G20").Value = 3.681611005342314
$ws.Range("I20").Value = 32.00217522400011
$ws.Range("L20").Value = 10.64155142968137
$ws.Range("M20").Value = 16.79531286797156
$ws.Range("B21").Value = 19.76742362987817
$ws.Range("C21").Value = 12.34128298482096
$ws.Range("D21").Value = 7.356905186639271
$ws.Range("F21").Value = 41.62239405521003
$ws.Range("G21").Value = 3.676901870259555
$ws.Range("I21").Value = 32.20489576235068
$ws.Range("L21").Value = 10.63537962422531
$ws.Range("M21").Value = 16.89201961858063
$ws.Range("B22").Value = 20.09672377231938
$ws.Range("C22").Value = 12.76895565916211
$ws.Range("D22").Value = 7.349712012626306
$ws.Range("F22").Value = 42.02097683356486
$ws.Range("G22").Value = 3.673933360156718
$ws.Range("I22").Value = 32.3426873056203
$ws.Range("L22").Value = 10.63262905676278
$ws.Range("M22").Value = 16.95828033396159
$ws.Range("B23").Value = 19.92120520677062
$ws.Range("C23").Value = 12.5425370740029
$ws.Range("D23").Value = 7.353489055842077
$ws.Range("F23").Value = 41.80784358014353
$ws.Range("G23").Value = 3.675507760696362
$ws.Range("I23").Value = 32.26866641659741
$ws.Range("L23").Value = 10.63397973944317
$ws.Range("M23").Value = 16.92264051247647
$ws.Range("B24").Value = 19.25198344192666
$ws.Range("C24").Value = 11.64453918135899
$ws.Range("D24").Value = 7.369152745927399
$ws.Range("F24").Value = 41.01007795374817
$ws.Range("G24").Value = 3.681689826094985
$ws.Range("I24").Value = 31.99895270031762
$ws.Range("L24").Value = 10.64167397570172
$ws.Range("M24").Value = 16.79378484196156
$ws.Range("B25").Value = 18.52517041535621
$ws.Range("C25").Value = 10.58890540389492
$ws.Range("D25").Value = 7.388917821900494
$ws.Range("F25").Value = 40.17357242651445
$ws.Range("G25").Value = 3.688831107098109
$ws.Range("I25").Value = 31.7313589555502
$ws.Range("L25").Value = 10.65548580079948
$ws.Range("M25").Value = 16.66834935743386
